# Regenerate merged AHB files
# 1) Rename the header columns: "..._old" -> "..._FV2304", "..._new" -> "..._FV2310"
# 2) Turn the sheet's used range into an Excel Table ("Table1")
# 3) Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header cells (row 1) -------------------------------------
$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)

$newHeaders = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

$fv2304Headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

$fv2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt 10; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $fv2304Headers[$i]
}

# Column 11 (K) is "diff" and stays unchanged.

for ($i = 0; $i -lt 10; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $fv2310Headers[$i]
}

# --- 2) Convert A1:U91 into a native Excel Table -------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U91"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3) Freeze the header row --------------------------------------------
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
